# Add a "TFP" parameter sheet to the production function workbook.
# Mirrors the structure of elasPROD (the other production-function elasticity
# sheet) but with every elasticity value set to 1 instead of 0.95, and makes
# the new sheet the active one.

$wb = $excel.ActiveWorkbook

# Capture the row labels from elasPROD (A2:A36) up front, as plain strings,
# before any sheet-collection mutation (Add/Move renumbers the sheets, which
# can shift what an already-bound worksheet reference resolves to).
$source = $wb.Worksheets.Item("elasPROD")
$labels = @()
for ($r = 2; $r -le 36; $r++) {
    $labels += $source.Cells.Item($r, 1).Text
}

# Add the new TFP sheet
$new = $wb.Worksheets.Add()
$new.Name = "TFP"

# Header label in B1
$new.Range("B1").Value = "TFP"

# Row labels (A2:A36) copied verbatim from elasPROD, values (B2:B36) set to 1.
# A leading apostrophe marks the label as explicit text (matching the
# quote-prefixed style elasPROD itself uses for these code-like labels).
for ($r = 2; $r -le 36; $r++) {
    $new.Cells.Item($r, 1).Value = "'" + $labels[$r - 2]
    $new.Cells.Item($r, 2).Value = 1
}

# Move the new sheet to the end of the tab strip (after elasPROD), then
# re-fetch the sheet collection fresh and make TFP the active/selected tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Move($null, $lastSheet)

$tfp = $wb.Worksheets.Item("TFP")
$tfp.Activate()
